# Fix glitch in log
# The "Random chance" column (D) for rows 76:90 was a shared formula (=1/3)
# left over from an earlier experiment; it should instead be the constant
# 0.25 that applies to these later test groups. Also backfill the missing
# "Date" (column G) entries for the last test group (rows 88:90), and move
# the active selection to G90 to reflect where the log now ends.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Log")

# Replace the shared "=1/3" formula in D76:D90 with the literal value 0.25.
$ws.Range("D76:D90").Value = 0.25

# Backfill the Date column for the last block of rows, which was missing it.
$ws.Range("G88:G90").Value = 42562

# Move the active cell/selection to reflect the new last logged row.
$ws.Range("G90").Select() | Out-Null
